$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("WIP")
$dst = $wb.Worksheets.Item("To Check")

$vals = $src.Range("B1117:Q1119").Value()
$src.Range("1117:1119").Delete()
$dst.Range("B1188:Q1190").Value = $vals
$dst.Range("B1188:B1190").NumberFormat = "@"
$dst.Range("D1188:D1190").NumberFormat = "@"

Write-Host ("UsedRange: " + $src.UsedRange.Address())
Write-Host ("AutoFilter before: " + $src.AutoFilter.Range.Address())
$src.Range("B1:R1285").AutoFilter()
Write-Host ("AutoFilter mid-null? " + ($src.AutoFilter -eq $null))
$src.Range("B1:R1285").AutoFilter()
if ($src.AutoFilter -ne $null) {
    Write-Host ("AutoFilter after: " + $src.AutoFilter.Range.Address())
}
